# Apply updated "想去人数" (F) / "最低票价" (G) figures scraped at commit 456a3b4
# for the gh-pages generated workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 539
$ws1.Range("F5").Value  = 6251
$ws1.Range("F6").Value  = 700
$ws1.Range("F8").Value  = 55
$ws1.Range("F9").Value  = 106
$ws1.Range("F10").Value = 297
$ws1.Range("F12").Value = 651
$ws1.Range("F13").Value = 1114
$ws1.Range("G13").Value = 61.2
$ws1.Range("F15").Value = 375
$ws1.Range("F17").Value = 10
$ws1.Range("F18").Value = 1400
$ws1.Range("F19").Value = 639
$ws1.Range("F20").Value = 5
$ws1.Range("F21").Value = 379
$ws1.Range("F23").Value = 1052
$ws1.Range("F24").Value = 104
$ws1.Range("F25").Value = 2148
$ws1.Range("F26").Value = 229
$ws1.Range("F27").Value = 73
$ws1.Range("F28").Value = 380
$ws1.Range("F29").Value = 53
$ws1.Range("F30").Value = 3469

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 28
$ws2.Range("F9").Value  = 691
$ws2.Range("F19").Value = 371
$ws2.Range("F20").Value = 310
$ws2.Range("F25").Value = 173
$ws2.Range("F33").Value = 1526

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value  = 1811
$ws3.Range("F6").Value  = 1175
$ws3.Range("F9").Value  = 426
$ws3.Range("F12").Value = 726

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1811
$ws4.Range("F4").Value  = 1175
$ws4.Range("F6").Value  = 426
$ws4.Range("F9").Value  = 726
$ws4.Range("F10").Value = 28
$ws4.Range("F11").Value = 539
$ws4.Range("F16").Value = 6251
$ws4.Range("F20").Value = 692
$ws4.Range("F21").Value = 106
$ws4.Range("F22").Value = 297
$ws4.Range("F24").Value = 651
$ws4.Range("F28").Value = 1114
$ws4.Range("G28").Value = 61.2
$ws4.Range("F29").Value = 375
$ws4.Range("F30").Value = 371
$ws4.Range("F31").Value = 310
$ws4.Range("F33").Value = 10
$ws4.Range("F34").Value = 1400
$ws4.Range("F35").Value = 639
$ws4.Range("F36").Value = 5
$ws4.Range("F37").Value = 379
$ws4.Range("F38").Value = 173
$ws4.Range("F42").Value = 2148
$ws4.Range("F46").Value = 1526
$ws4.Range("F47").Value = 229
$ws4.Range("F48").Value = 380
$ws4.Range("F49").Value = 53
$ws4.Range("F50").Value = 3469
